# Translation fixes (fixes #2060, #2067, #2073, #2074, #2076, 2077, #2078, #2079, #2080)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the German translation for "description_comment" (row 40, col B)
$ws.Range("B40").Value = "Kommentar"

# Shorten the German translation for "project_affiliation" (row 41, col B)
$ws.Range("B41").Value = "Projekt"

# Add a new translation key "project" with German value "Projekt" as a new row
$ws.Range("A62").Value = "project"
$ws.Range("B62").Value = "Projekt"

# Update the view: unfreeze the header rows and move the selection/scroll position
try {
    $aw = $excel.ActiveWindow
    $aw.FreezePanes = $false
    $aw.ScrollRow = 35
    $aw.ScrollColumn = 1
    [void]$ws.Range("B46").Select()
} catch {
}
